$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 337, shifting existing rows 337:362 down to 338:363.
$ws.Rows(337).Insert()

# Populate the newly inserted row 337 with its data.
$ws.Range("A337").Value = 9
$ws.Range("B337").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C337").Value = "Metropolitana"
$ws.Range("D337").Value = 44918
$ws.Range("E337").Value = 13
$ws.Range("F337").Value = 100112021
$ws.Range("G337").Value = "Ají"
$ws.Range("H337").Value = "Americana (o)"
$ws.Range("I337").Value = "Primera"
$ws.Range("J337").Value = 70
$ws.Range("K337").Value = 34000
$ws.Range("L337").Value = 36000
$ws.Range("M337").Value = 35000
$ws.Range("N337").Value = "`$/caja 25 kilos"
$ws.Range("O337").Value = "Provincia de Limarí"
$ws.Range("P337").Value = 1400
$ws.Range("Q337").Value = 25
$ws.Range("R337").Value = "Hortaliza"
